$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Replace the old "Przykladowy przebieg formatowania tekstu:" paragraph
#    with the new "Implementacja aplikacji" block (library description,
#    function list) followed by a blank paragraph and the (re-worded)
#    "Przykladowy przebieg procedury formatowania tekstu:" paragraph.
# ---------------------------------------------------------------------------
$newBlock = "Implementacja aplikacji: ^p" + `
    "^tBiblioteka aplikacji składa się z plików:^p" + `
    "^tformat.h – deklaracje funkcji formatujących^p" + `
    "^tformat.c – implementacje funkcji formatujących^p" + `
    "W skład aplikacji wchodzą funkcje:^p" + `
    "^ttoLowCases^t- konwersja tekstu na małe litery^p" + `
    "^tformatCapitalLetter – wielkie litery na początku zdań^p" + `
    "^tformatRepeatedLetters – usunięcie powtarzających się znaków^p" + `
    "^tformatSpaces – formatowanie spacji^p" + `
    "^tremoveBlankLetters – finalne usunięcie z tablicy nadmiarowych znaków^p" + `
    "^p" + `
    "Przykładowy przebieg procedury formatowania tekstu:"

$target = $d.Paragraphs.Item(10).Range
$target.Find.Execute("Przykładowy przebieg formatowania tekstu:", $true, $false, $false, $false, $false, `
    $true, 1, $false, $newBlock, 2) | Out-Null

# ---------------------------------------------------------------------------
# 2) Move the "_GoBack" bookmark from the old "Tekst wyjsciowy:" paragraph to
#    the new "Przykladowy przebieg procedury formatowania tekstu:" paragraph,
#    right after the word "procedury".
# ---------------------------------------------------------------------------
$d.Bookmarks.Item("_GoBack").Delete()

$procRange = $d.Content
$procRange.Find.ClearFormatting()
$procRange.Find.Execute("Przykładowy przebieg procedury") | Out-Null
$bmPos = $procRange.End
$bmRange = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null

# ---------------------------------------------------------------------------
# 3) Replace the trailing empty paragraph at the end of the document with
#    the final "Tekst wyjsciowy:" and quote paragraphs.
# ---------------------------------------------------------------------------
$lastIdx = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($lastIdx)
$lastPara.Range.InsertParagraphBefore()

$outLabelPara = $d.Paragraphs.Item($lastIdx)
$outLabelPara.Range.Text = "Tekst wyjściowy: "

$quotePara = $d.Paragraphs.Item($lastIdx + 1)
$quotePara.Range.Text = "`"Ala ma kota, 1123. Mama ma psa.  `""
$quotePara.Range.ParagraphFormat.FirstLineIndent = 35.4

Write-Output "Done. Paragraph count: $($d.Paragraphs.Count)"
